$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the question/answer text for the gender-roles row (V202287)
$ws.Range("B6").Value = "How would this respondent assess whether it is better, worse, or makes no difference for the family as a whole if the man works outside the home and the woman takes care of the home and family?"
$ws.Range("C6").Value = "1. Better 2. Worse 3. Makes no difference"

# Update the active cell selection to match the saved view state
$ws.Range("G20").Select()
